# "Generate Report for Handback"
#
# The localization-status report is regenerated: the zh-cn and de-de
# handback rows move from "Ready for handoff" to "Handed back: in sync
# with en-US", their "Latest Handback DateTime" stamps advance to the
# new handback run, and the stale "handback file is not the latest"
# error detail is cleared now that the handback is in sync. The
# Overview roll-up sheet mirrors the same status text for both locales.
# Status/Error-Detail columns are also widened/narrowed to fit the new
# text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$overview = $wb.Sheets("Overview")
$zhcn     = $wb.Sheets("zh-cn")
$dede     = $wb.Sheets("de-de")

# --- zh-cn sheet --------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-28 22:49:15"
$zhcn.Range("P2").Value = ""

# --- de-de sheet ---------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-28 22:49:21"
$dede.Range("P2").Value = ""

# --- Overview roll-up sheet (zh-cn / de-de status columns) ---------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- Column width adjustments --------------------------------------------
# Status columns widen, Error Detail columns narrow to fit the new text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
